$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (shifts old rows 13-23 down to 14-24)
$ws.Rows.Item(13).Insert()

# Row 13 (new): clear the stray carried-over cell in col A, then set B13/C13
$ws.Range("A13").Clear()
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C13").PasteSpecial(-4122) | Out-Null
$ws.Range("B13").Value = '5817650 - Érica Leonor Romão'
$ws.Range("C13").Value = '5817650 - Érica Leonor Romão'

# Fix Objetivos (row 10): was showing the Erica value, now the correct PT objectives text
$ws.Range("B10").Value = 'Proporcionar aos alunos uma visão prática do licenciamento ambiental verificando sua interface com os demais instrumentos da política ambiental. Transmitir aos alunos o ferramental teórico e prático necessário para sua atuação profissional.'
$ws.Range("C10").Value = 'Proporcionar aos alunos uma visão prática do licenciamento ambiental verificando sua interface com os demais instrumentos da política ambiental. Transmitir aos alunos o ferramental teórico e prático necessário para sua atuação profissional.'

# Fix Programa resumido (row 14): was showing "Semestral", now the correct summary text
$ws.Range("B14").Value = 'O licenciamento ambiental como instrumento da política ambiental; procedimentos do licenciamento ambiental; licenciamento no estado de São Paulo: aspectos institucionais e legislação aplicada; estudo de caso durante as etapas do licenciamento.'
$ws.Range("C14").Value = 'O licenciamento ambiental como instrumento da política ambiental; procedimentos do licenciamento ambiental; licenciamento no estado de São Paulo: aspectos institucionais e legislação aplicada; estudo de caso durante as etapas do licenciamento.'

# Fix Programa (row 16): was showing a date, now the correct programme text
$ws.Range("B16").Value = 'O licenciamento ambiental como instrumento da política nacional do meio ambiente; etapas e prazos do licenciamento ambiental: licença prévia, de instalação e de operação do empreendimento; empreendimentos sujeitos ao licenciamento ambiental; o licenciamento no estado de São Paulo: aspectos institucionais, legislação aplicada e documentação; abordagem técnica e legal no âmbito do licenciamento ambiental referente à vegetação nativa e área de preservação permanente no Estado de São Paulo, estudo de caso; atuação do engenheiro ambiental.'
$ws.Range("C16").Value = 'O licenciamento ambiental como instrumento da política nacional do meio ambiente; etapas e prazos do licenciamento ambiental: licença prévia, de instalação e de operação do empreendimento; empreendimentos sujeitos ao licenciamento ambiental; o licenciamento no estado de São Paulo: aspectos institucionais, legislação aplicada e documentação; abordagem técnica e legal no âmbito do licenciamento ambiental referente à vegetação nativa e área de preservação permanente no Estado de São Paulo, estudo de caso; atuação do engenheiro ambiental.'

# Fix Método (row 19): was showing the Erica value, now the correct method text
$ws.Range("B19").Value = 'Aulas teóricas e práticas, trabalhos de campo, exercícios dirigidos e seminários.As avaliações serão por meio de trabalhos em equipe, provas individuais conforme adequação ao conteúdo.'
$ws.Range("C19").Value = 'Aulas teóricas e práticas, trabalhos de campo, exercícios dirigidos e seminários.As avaliações serão por meio de trabalhos em equipe, provas individuais conforme adequação ao conteúdo.'

# Fix Critério (row 20): now the weighted-average text
$ws.Range("B20").Value = 'Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios.'
$ws.Range("C20").Value = 'Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios.'

# Fix Norma de recuperação (row 21): now the recuperation-grade text
$ws.Range("B21").Value = 'A nota final será composta pela média obtida da nota do período somada à nota de recuperação'
$ws.Range("C21").Value = 'A nota final será composta pela média obtida da nota do período somada à nota de recuperação'

# Fix Bibliografia (row 22): now the full bibliography text
$ws.Range("B22").Value = 'Bibliografia básica:SÁNCHEZ, L.E., Avaliação de impacto ambiental: conceitos e métodos, Ed. Oficina de textos, 3° reimpressão, 2011.OLIVEIRA, A.I.A., O licenciamento ambiental, Iglu Editora, 1999.BECHARA, E. Licenciamento e Compensação Ambiental – Ed. Atlas, 2009.CURI, D. (Org.), Gestão ambiental, Ed. Pearson, 2011.LIMA, A., Zoneamento Ecológico-Econômico – a luz dos direitos socioambientais, Juruá Editora, 2006.MOURA, L.A.A., Qualidade e Gestão ambiental – sustentabilidade e ISO 14.001, 6° ed., Ed. Del Rey, 2011. SEIFFERT, M.E.B., Gestão ambiental: instrumentos, esferas de ação e educação ambiental, Editora Atlas, 2007.Bibliografia complementar:BRAGA B. (Org.), Introdução à engenharia ambiental: o desafio do desenvolvimento sustentável, 2° ed., Ed. Pearson Prentice Hall, 2005CALIJURI, M.C., CUNHA, D.G.F. (Org.), Engenharia ambiental: conceitos, tecnologia e gestão, Ed. Campus, 2013KRAWULSKI, C.C., FEIJÓ, C.C.C., Introdução à gestão ambiental, Ed. Pearson, 2009CETESB  Manuais de licenciamento ambiental'
$ws.Range("C22").Value = 'Bibliografia básica:SÁNCHEZ, L.E., Avaliação de impacto ambiental: conceitos e métodos, Ed. Oficina de textos, 3° reimpressão, 2011.OLIVEIRA, A.I.A., O licenciamento ambiental, Iglu Editora, 1999.BECHARA, E. Licenciamento e Compensação Ambiental – Ed. Atlas, 2009.CURI, D. (Org.), Gestão ambiental, Ed. Pearson, 2011.LIMA, A., Zoneamento Ecológico-Econômico – a luz dos direitos socioambientais, Juruá Editora, 2006.MOURA, L.A.A., Qualidade e Gestão ambiental – sustentabilidade e ISO 14.001, 6° ed., Ed. Del Rey, 2011. SEIFFERT, M.E.B., Gestão ambiental: instrumentos, esferas de ação e educação ambiental, Editora Atlas, 2007.Bibliografia complementar:BRAGA B. (Org.), Introdução à engenharia ambiental: o desafio do desenvolvimento sustentável, 2° ed., Ed. Pearson Prentice Hall, 2005CALIJURI, M.C., CUNHA, D.G.F. (Org.), Engenharia ambiental: conceitos, tecnologia e gestão, Ed. Campus, 2013KRAWULSKI, C.C., FEIJÓ, C.C.C., Introdução à gestão ambiental, Ed. Pearson, 2009CETESB  Manuais de licenciamento ambiental'

